$wb = $excel.ActiveWorkbook

# Rename sheet "iterationdata" -> "iterationData"
$wsData = $wb.Worksheets.Item("iterationdata")
$wsData.Name = "iterationData"

$wsRun = $wb.Worksheets.Item("runManager")

# Update cell values B3 and B6 on iterationData sheet: "yes" -> "no"
$wsData.Range("B3").Value = "no"
$wsData.Range("B6").Value = "no"

# Change active/selected sheet and selection state
# Previously iterationData was the tab-selected sheet with selection B4.
# Now runManager should be tab-selected/active with selection B1.
$wsData.Activate()
$wsData.Range("A1").Select()
$wsRun.Activate()
$wsRun.Range("B1").Select()
